# Update NATMI LR-pair TPM-derived statistics for Spon2-Itga4 (rows 2-17)
# per commit "update scripts wuth new tpm"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.150782
$ws.Range("H2").Value = 3.452345999999999
$ws.Range("I2").Value = 0.03823856951930295
$ws.Range("J2").Value = 0.03823856951930295
$ws.Range("M2").Value = 0.1419263333333333
$ws.Range("N2").Value = 0.425779
$ws.Range("O2").Value = 0.002583058778296354
$ws.Range("P2").Value = 0.002583058778296354
$ws.Range("Q2").Value = 0.163326269726
$ws.Range("R2").Value = 1.469936427534
$ws.Range("S2").Value = 0.00009877247266633089
$ws.Range("T2").Value = 0.00009877247266633086

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.150782
$ws.Range("H3").Value = 3.452345999999999
$ws.Range("I3").Value = 0.03823856951930295
$ws.Range("J3").Value = 0.03823856951930295
$ws.Range("O3").Value = 0.001399682868699959
$ws.Range("P3").Value = 0.001399682868699959
$ws.Range("Q3").Value = 0.08850165689799998
$ws.Range("R3").Value = 0.7965149120819999
$ws.Range("S3").Value = 0.00005352187067976077
$ws.Range("T3").Value = 0.00005352187067976076

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.150782
$ws.Range("H4").Value = 3.452345999999999
$ws.Range("I4").Value = 0.03823856951930295
$ws.Range("J4").Value = 0.03823856951930295
$ws.Range("M4").Value = 2.613991
$ws.Range("N4").Value = 7.841973
$ws.Range("O4").Value = 0.04757462720522382
$ws.Range("P4").Value = 0.04757462720522382
$ws.Range("Q4").Value = 3.008133790961999
$ws.Range("R4").Value = 27.073204118658
$ws.Range("S4").Value = 0.001819185689741873
$ws.Range("T4").Value = 0.001819185689741872

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.150782
$ws.Range("H5").Value = 3.452345999999999
$ws.Range("I5").Value = 0.03823856951930295
$ws.Range("J5").Value = 0.03823856951930295
$ws.Range("M5").Value = 52.11224233333333
$ws.Range("N5").Value = 156.336727
$ws.Range("O5").Value = 0.9484426311477799
$ws.Range("P5").Value = 0.9484426311477798
$ws.Range("Q5").Value = 59.96983045683799
$ws.Range("R5").Value = 539.7284741115419
$ws.Range("S5").Value = 0.03626708948621499
$ws.Range("T5").Value = 0.03626708948621498

# Row 6
$ws.Range("I6").Value = 0.9169230158851821
$ws.Range("J6").Value = 0.916923015885182
$ws.Range("M6").Value = 0.1419263333333333
$ws.Range("N6").Value = 0.425779
$ws.Range("O6").Value = 0.002583058778296354
$ws.Range("P6").Value = 0.002583058778296354
$ws.Range("Q6").Value = 3.916402148224778
$ws.Range("R6").Value = 35.24761933402301
$ws.Range("S6").Value = 0.002368466045204187
$ws.Range("T6").Value = 0.002368466045204186

# Row 7
$ws.Range("I7").Value = 0.9169230158851821
$ws.Range("J7").Value = 0.916923015885182
$ws.Range("O7").Value = 0.001399682868699959
$ws.Range("P7").Value = 0.001399682868699959
$ws.Range("S7").Value = 0.00128340143725119
$ws.Range("T7").Value = 0.001283401437251189

# Row 8
$ws.Range("I8").Value = 0.9169230158851821
$ws.Range("J8").Value = 0.916923015885182
$ws.Range("M8").Value = 2.613991
$ws.Range("N8").Value = 7.841973
$ws.Range("O8").Value = 0.04757462720522382
$ws.Range("P8").Value = 0.04757462720522382
$ws.Range("Q8").Value = 72.13206828782234
$ws.Range("R8").Value = 649.188614590401
$ws.Range("S8").Value = 0.04362227065662706
$ws.Range("T8").Value = 0.04362227065662705

# Row 9
$ws.Range("I9").Value = 0.9169230158851821
$ws.Range("J9").Value = 0.916923015885182
$ws.Range("M9").Value = 52.11224233333333
$ws.Range("N9").Value = 156.336727
$ws.Range("O9").Value = 0.9484426311477799
$ws.Range("P9").Value = 0.9484426311477798
$ws.Range("Q9").Value = 1438.017125009055
$ws.Range("R9").Value = 12942.1541250815
$ws.Range("S9").Value = 0.8696488777460998
$ws.Range("T9").Value = 0.8696488777460996

# Row 10
$ws.Range("G10").Value = 1.290098666666667
$ws.Range("H10").Value = 3.870296
$ws.Range("I10").Value = 0.04286783035543951
$ws.Range("J10").Value = 0.0428678303554395
$ws.Range("M10").Value = 0.1419263333333333
$ws.Range("N10").Value = 0.425779
$ws.Range("O10").Value = 0.002583058778296354
$ws.Range("P10").Value = 0.002583058778296354
$ws.Range("Q10").Value = 0.1830989733982222
$ws.Range("R10").Value = 1.647890760584
$ws.Range("S10").Value = 0.0001107301255061369
$ws.Range("T10").Value = 0.0001107301255061369

# Row 11
$ws.Range("G11").Value = 1.290098666666667
$ws.Range("H11").Value = 3.870296
$ws.Range("I11").Value = 0.04286783035543951
$ws.Range("J11").Value = 0.0428678303554395
$ws.Range("O11").Value = 0.001399682868699959
$ws.Range("P11").Value = 0.001399682868699959
$ws.Range("Q11").Value = 0.09921589802577777
$ws.Range("R11").Value = 0.8929430822319999
$ws.Range("S11").Value = 0.00006000136776684475
$ws.Range("T11").Value = 0.00006000136776684474

# Row 12
$ws.Range("G12").Value = 1.290098666666667
$ws.Range("H12").Value = 3.870296
$ws.Range("I12").Value = 0.04286783035543951
$ws.Range("J12").Value = 0.0428678303554395
$ws.Range("M12").Value = 2.613991
$ws.Range("N12").Value = 7.841973
$ws.Range("O12").Value = 0.04757462720522382
$ws.Range("P12").Value = 0.04757462720522382
$ws.Range("Q12").Value = 3.372306303778666
$ws.Range("R12").Value = 30.350756734008
$ws.Range("S12").Value = 0.002039421048256812
$ws.Range("T12").Value = 0.002039421048256811

# Row 13
$ws.Range("G13").Value = 1.290098666666667
$ws.Range("H13").Value = 3.870296
$ws.Range("I13").Value = 0.04286783035543951
$ws.Range("J13").Value = 0.0428678303554395
$ws.Range("M13").Value = 52.11224233333333
$ws.Range("N13").Value = 156.336727
$ws.Range("O13").Value = 0.9484426311477799
$ws.Range("P13").Value = 0.9484426311477798
$ws.Range("Q13").Value = 67.22993435124354
$ws.Range("R13").Value = 605.069409161192
$ws.Range("S13").Value = 0.04065767781390971
$ws.Range("T13").Value = 0.04065767781390971

# Row 14
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.05930433333333333
$ws.Range("H14").Value = 0.177913
$ws.Range("I14").Value = 0.001970584240075516
$ws.Range("J14").Value = 0.001970584240075516
$ws.Range("M14").Value = 0.1419263333333333
$ws.Range("N14").Value = 0.425779
$ws.Range("O14").Value = 0.002583058778296354
$ws.Range("P14").Value = 0.002583058778296354
$ws.Range("Q14").Value = 0.008416846580777777
$ws.Range("R14").Value = 0.075751619227
$ws.Range("S14").Value = 0.000005090134919699512
$ws.Range("T14").Value = 0.000005090134919699509

# Row 15
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.05930433333333333
$ws.Range("H15").Value = 0.177913
$ws.Range("I15").Value = 0.001970584240075516
$ws.Range("J15").Value = 0.001970584240075516
$ws.Range("O15").Value = 0.001399682868699959
$ws.Range("P15").Value = 0.001399682868699959
$ws.Range("Q15").Value = 0.004560839291222221
$ws.Range("R15").Value = 0.04104755362099999
$ws.Range("S15").Value = 0.000002758193002163827
$ws.Range("T15").Value = 0.000002758193002163826

# Row 16
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.05930433333333333
$ws.Range("H16").Value = 0.177913
$ws.Range("I16").Value = 0.001970584240075516
$ws.Range("J16").Value = 0.001970584240075516
$ws.Range("M16").Value = 2.613991
$ws.Range("N16").Value = 7.841973
$ws.Range("O16").Value = 0.04757462720522382
$ws.Range("P16").Value = 0.04757462720522382
$ws.Range("Q16").Value = 0.1550209935943333
$ws.Range("R16").Value = 1.395188942349
$ws.Range("S16").Value = 0.00009374981059808194
$ws.Range("T16").Value = 0.00009374981059808193

# Row 17
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.05930433333333333
$ws.Range("H17").Value = 0.177913
$ws.Range("I17").Value = 0.001970584240075516
$ws.Range("J17").Value = 0.001970584240075516
$ws.Range("M17").Value = 52.11224233333333
$ws.Range("N17").Value = 156.336727
$ws.Range("O17").Value = 0.9484426311477799
$ws.Range("P17").Value = 0.9484426311477798
$ws.Range("Q17").Value = 3.090481790083444
$ws.Range("R17").Value = 27.814336110751
$ws.Range("S17").Value = 0.001868986101555571
$ws.Range("T17").Value = 0.00186898610155557

